# Auto-generated edit script applying the Zeromus_Profits.xlsx diff
# Updates leve-profit calculation columns (H..N) across 8 profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3484.625
$ws.Range("I76").Value = 3038.889
$ws.Range("J76").Value = 4057.7144
$ws.Range("K76").Value = 3038.889
$ws.Range("L76").Value = 4057.7144
$ws.Range("M76").Value = -2723.889
$ws.Range("N76").Value = -4687.7144

$ws.Range("H79").Value = 3484.625
$ws.Range("I79").Value = 3038.889
$ws.Range("J79").Value = 4057.7144
$ws.Range("K79").Value = 3038.889
$ws.Range("L79").Value = 4057.7144
$ws.Range("M79").Value = -1946.889
$ws.Range("N79").Value = -6241.7144

$ws.Range("H80").Value = 1230.4615
$ws.Range("I80").Value = 742.56525
$ws.Range("J80").Value = 1931.8125
$ws.Range("K80").Value = 2227.69575
$ws.Range("L80").Value = 5795.4375
$ws.Range("M80").Value = -1229.69575
$ws.Range("N80").Value = -7791.4375

$ws.Range("H83").Value = 1230.4615
$ws.Range("I83").Value = 742.56525
$ws.Range("J83").Value = 1931.8125
$ws.Range("K83").Value = 6683.08725
$ws.Range("L83").Value = 17386.3125
$ws.Range("M83").Value = -1691.08725
$ws.Range("N83").Value = -27370.3125

$ws.Range("H86").Value = 7033.5
$ws.Range("I86").Value = 11485.6
$ws.Range("J86").Value = 2581.4
$ws.Range("K86").Value = 11485.6
$ws.Range("L86").Value = 2581.4
$ws.Range("M86").Value = -10362.6
$ws.Range("N86").Value = -4827.4

$ws.Range("H89").Value = 7033.5
$ws.Range("I89").Value = 11485.6
$ws.Range("J89").Value = 2581.4
$ws.Range("K89").Value = 57428
$ws.Range("L89").Value = 12907
$ws.Range("M89").Value = -51812
$ws.Range("N89").Value = -24139

$ws.Range("H138").Value = 224782.84
$ws.Range("I138").Value = 477009.47
$ws.Range("J138").Value = 4084.5417
$ws.Range("K138").Value = 1431028.41
$ws.Range("L138").Value = 12253.6251
$ws.Range("M138").Value = -1425888.41
$ws.Range("N138").Value = -22533.6251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H110").Value = 10068.4375
$ws.Range("I110").Value = 13083.272
$ws.Range("K110").Value = 13083.272
$ws.Range("M110").Value = -11038.272

$ws.Range("H122").Value = 2724.625
$ws.Range("I122").Value = 2391.8462
$ws.Range("J122").Value = 4166.6665
$ws.Range("K122").Value = 7175.5386
$ws.Range("L122").Value = 12499.9995
$ws.Range("M122").Value = -4725.5386
$ws.Range("N122").Value = -17399.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1275.5
$ws.Range("I99").Value = 1247.4375
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1247.4375
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 250.5625
$ws.Range("N99").Value = -4496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4468191.5
$ws.Range("I31").Value = 7445480.5
$ws.Range("J31").Value = 2258.611
$ws.Range("K31").Value = 7445480.5
$ws.Range("L31").Value = 2258.611
$ws.Range("M31").Value = -7445185.5
$ws.Range("N31").Value = -2848.611

$ws.Range("H34").Value = 4468191.5
$ws.Range("I34").Value = 7445480.5
$ws.Range("J34").Value = 2258.611
$ws.Range("K34").Value = 7445480.5
$ws.Range("L34").Value = 2258.611
$ws.Range("M34").Value = -7445278.5
$ws.Range("N34").Value = -2662.611

$ws.Range("H58").Value = 1835.091
$ws.Range("I58").Value = 1085.625
$ws.Range("J58").Value = 2263.3572
$ws.Range("K58").Value = 1085.625
$ws.Range("L58").Value = 2263.3572
$ws.Range("M58").Value = -882.625
$ws.Range("N58").Value = -2669.3572

$ws.Range("H136").Value = 1835.091
$ws.Range("I136").Value = 1085.625
$ws.Range("J136").Value = 2263.3572
$ws.Range("K136").Value = 3256.875
$ws.Range("L136").Value = 6790.071599999999
$ws.Range("M136").Value = -706.875
$ws.Range("N136").Value = -11890.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 450
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -3996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws.Range("H122").Value = 4253.857
$ws.Range("I122").Value = 3729
$ws.Range("J122").Value = 5065
$ws.Range("K122").Value = 11187
$ws.Range("L122").Value = 15195
$ws.Range("M122").Value = -8737
$ws.Range("N122").Value = -20095

$ws.Range("H132").Value = 2438.28
$ws.Range("I132").Value = 1943
$ws.Range("J132").Value = 3490.75
$ws.Range("K132").Value = 5829
$ws.Range("L132").Value = 10472.25
$ws.Range("M132").Value = -3299
$ws.Range("N132").Value = -15532.25

$ws.Range("H134").Value = 52660
$ws.Range("J134").Value = 52660
$ws.Range("L134").Value = 157980
$ws.Range("N134").Value = -163050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2080
$ws.Range("I7").Value = 1400
$ws.Range("J7").Value = 2250
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 2250
$ws.Range("M7").Value = -1288
$ws.Range("N7").Value = -2474

$ws.Range("H46").Value = 797.6667
$ws.Range("I46").Value = 797
$ws.Range("J46").Value = 797.8
$ws.Range("K46").Value = 797
$ws.Range("L46").Value = 797.8
$ws.Range("M46").Value = -609
$ws.Range("N46").Value = -1173.8

$ws.Range("H55").Value = 313.66666
$ws.Range("I55").Value = 370.5
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 370.5
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = -197.5
$ws.Range("N55").Value = -546

$ws.Range("H92").Value = 29800
$ws.Range("J92").Value = 29800
$ws.Range("L92").Value = 29800
$ws.Range("N92").Value = -34792

$ws.Range("H126").Value = 2080
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -11690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 58825360
$ws.Range("I122").Value = 76924830
$ws.Range("J122").Value = 2077.5
$ws.Range("K122").Value = 230774490
$ws.Range("L122").Value = 6232.5
$ws.Range("M122").Value = -230772040
$ws.Range("N122").Value = -11132.5

$ws.Range("H132").Value = 1516.674
$ws.Range("I132").Value = 1082.0968
$ws.Range("J132").Value = 2414.8
$ws.Range("K132").Value = 3246.2904
$ws.Range("L132").Value = 7244.400000000001
$ws.Range("M132").Value = -716.2903999999999
$ws.Range("N132").Value = -12304.4
